# Add the English version of the newly-added "Julys" (water drop company /
# medical social worker) resume entry to the english_version sheet.
# The same row already exists (in Chinese) at row 4 of the "main" sheet;
# here we insert a new row 4 in "english_version" and copy that row's
# values across, then update view-state / defined names to match.

$wb = $excel.ActiveWorkbook
$main = $wb.Worksheets.Item("main")
$ws = $wb.Worksheets.Item("english_version")

# Insert a new blank row above the current row 4, shifting existing
# rows 4-21 down to 5-22.
$ws.Rows.Item(4).Insert()

# Copy the corresponding row (A4:J4) from "main" into the new row.
$srcRange = $main.Range("A4:J4")
$srcRange.Copy()
$dstRange = $ws.Range("A4:J4")
$dstRange.PasteSpecial()

# The entry that used to be row 4 (now shifted to row 5) was still marked
# as ongoing ("end" = 9999) in english_version; "main" already records
# that it ended in 2022 once the new job started, so bring the two in
# sync here too.
$ws.Range("G5").Value = $main.Range("G5").Value2

# Update the two defined names that pointed into english_version so they
# keep referencing the same logical cells now that rows shifted down by one.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -eq "english_version!OLE_LINK52") {
        $n.RefersTo = "=english_version!`$C`$8"
    }
    if ($n.Name -eq "english_version!OLE_LINK65") {
        $n.RefersTo = "=english_version!`$D`$10"
    }
}

# Match the resulting selection / active-sheet state: the "main" sheet had
# its whole 4th row selected, and "english_version" becomes the active tab
# with B4 selected.
$main.Rows.Item(4).Select()
$ws.Select()
$ws.Range("B4").Select()
